# Remove duplicate galaxy-cut entries from both sheets, renumber the
# remaining IDs, and drop the now-unused picture anchors for the trailing
# (now-deleted) entries on "Sheet1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Sheet": one data row per entry (header row 1, data rows 2-28).
# Duplicate entries sit at data rows 7, 8, 15, 18, 19 (entries #6, #7,
# #14, #17, #18) - each duplicates an earlier entry (#1-#5).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet")

$ws1.Rows("19:19").Delete()
$ws1.Rows("18:18").Delete()
$ws1.Rows("15:15").Delete()
$ws1.Rows("8:8").Delete()
$ws1.Rows("7:7").Delete()

# Renumber the trailing "_N" id suffix in column A for every remaining
# data row (2-23) so the ids stay sequential (1..22).
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws1.Range("A$r")
    $txt = $cell.Text
    $lastUnderscore = $txt.LastIndexOf("_")
    $prefix = $txt.Substring(0, $lastUnderscore)
    $cell.Value = $prefix + "_" + ($r - 1)
}

# ---------------------------------------------------------------------
# Sheet "Sheet1": one labeled row every 13 rows (row 2, 15, 28, ...),
# each followed by picture anchors for 4 extra filter views. The same
# duplicate entries (#6, #7, #14, #17, #18) occupy the 13-row blocks
# starting at rows 67, 80, 171, 210, 223.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet1")

$ws2.Rows("223:235").Delete()
$ws2.Rows("210:222").Delete()
$ws2.Rows("171:183").Delete()
$ws2.Rows("80:92").Delete()
$ws2.Rows("67:79").Delete()

# Renumber the trailing "_N" id suffix in column A for every remaining
# labeled row so the ids stay sequential (1..22).
$labelRows = @(2, 15, 28, 41, 54, 67, 80, 93, 106, 119, 132, 145, 158, 171, 184, 197, 210, 223, 236, 249, 262, 275)
for ($i = 0; $i -lt $labelRows.Count; $i++) {
    $r = $labelRows[$i]
    $cell = $ws2.Range("A$r")
    $txt = $cell.Text
    $lastUnderscore = $txt.LastIndexOf("_")
    $prefix = $txt.Substring(0, $lastUnderscore)
    $cell.Value = $prefix + "_" + ($i + 1)
}

# Remove the picture anchors that belonged to the 5 entries that no
# longer exist (previously entries #23-#27, images 111-135).
for ($i = 111; $i -le 135; $i++) {
    $ws2.Shapes.Item("Image $i").Delete()
}
